$d = $word.ActiveDocument

# Locate the paragraph that ends the "Doctor prescription" list item so we
# can insert the new "Customer Feedback" list item directly after it.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "Doctor prescription") {
        $target = $para
    }
}

if ($target -ne $null) {
    # InsertParagraphAfter clones the paragraph's pPr/numPr/rPr (list style,
    # numbering, fonts, size) onto the new paragraph, matching the sibling
    # list items already in the document.
    $target.Range.InsertParagraphAfter()

    # The freshly inserted paragraph is now the one right after $target.
    $newPara = $target.Next()
    $newPara.Range.InsertAfter("Customer Feedback")
}
